$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "tool" (column AC) value for the existing "Schleiermacher digital" rows (77-79) ---
$ws.Range("AC77").Value = "ediarum.BASE.edit, ediarum.DB, ediarum.WEB"
$ws.Range("AC78").Value = "ediarum.BASE.edit, ediarum.DB, ediarum.WEB"
$ws.Range("AC79").Value = "ediarum.BASE.edit, ediarum.DB, ediarum.WEB"

# --- New row 91: Bolland / Cosijn correspondence edition ---
$ws.Range("A91").Value = 60
$ws.Range("B91").Value = "`nDe filosoof en de filoloog: De correspondentie tussen G. J. P. J. Bolland en P. J. Cosijn (1879-1899)"
$ws.Range("C91").Value = "Bolland"
$ws.Range("D91").Value = "Gerardus"
$ws.Range("E91").Value = "Porck"
$ws.Range("F91").Value = "Thijs"
$ws.Range("G91").Value = "Leiden University Fund (LUF)"
$ws.Range("H91").Value = 52.160114
$ws.Range("I91").Value = 4.49701
$ws.Range("J91").Value = "Leiden"
$ws.Range("K91").Value = 2751773
$ws.Range("L91").Value = 52.15833
$ws.Range("M91").Value = 4.49306
$ws.Range("N91").Value = "Netherlands"
$ws.Range("O91").Value = "https://www.luf.nl/en"
$ws.Range("P91").Value = 2019
$ws.Range("Q91").Value = '"NLD"'
$ws.Range("R91").Value = "https://correspondentie-bolland-en-cosijn.huygens.knaw.nl/"
$ws.Range("S91").Value = "Modern"
$ws.Range("T91").Value = "History of Philosophy"
$ws.Range("U91").Value = "Letters"
$ws.Range("V91").Value = "yes"
$ws.Range("W91").Value = "yes"
$ws.Range("X91").Value = "yes"
$ws.Range("Y91").Value = "no"
$ws.Range("Z91").Value = "not provided"
$ws.Range("AA91").Value = "ony pdf"
$ws.Range("AB91").Value = "not provided"

# --- New row 92: Ernst Haeckel letter edition ---
$ws.Range("A92").Value = 61
$ws.Range("B92").Value = "Ernst Haeckel (1834–1919): Briefedition"
$ws.Range("C92").Value = "Haeckel"
$ws.Range("D92").Value = "Ernst"
$ws.Range("E92").Value = "Bach"
$ws.Range("F92").Value = "Thomas"
$ws.Range("G92").Value = "Friedrich-Schiller-Universität Jena"
$ws.Range("H92").Value = 50.9271
$ws.Range("I92").Value = 11.5892
$ws.Range("J92").Value = "Jena"
$ws.Range("K92").Value = 2895044
$ws.Range("L92").Value = 50.92878
$ws.Range("M92").Value = 11.5899
$ws.Range("N92").Value = "Germany"
$ws.Range("O92").Value = "https://www.uni-jena.de/"
$ws.Range("P92").Value = 2017
$ws.Range("Q92").Value = '"GER"'
$ws.Range("R92").Value = "https://haeckel-briefwechsel-projekt.uni-jena.de/de"
$ws.Range("S92").Value = "Modern"
$ws.Range("T92").Value = "History of Science"
$ws.Range("U92").Value = "Letters"
$ws.Range("V92").Value = "yes"
$ws.Range("W92").Value = "yes"
$ws.Range("X92").Value = "yes"
$ws.Range("Y92").Value = "no"
$ws.Range("Z92").Value = "XML-TEI"
$ws.Range("AA92").Value = "yes"
$ws.Range("AB92").Value = "yes"

# --- New row 93: placeholder ID only ---
$ws.Range("A93").Value = 62

# The multi-line title in B91 makes the runtime auto-expand the row height;
# restore it to the sheet's standard (non-custom) height to match the source file.
$ws.Rows(91).AutoFit()
